$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph — this
# text, along with the blank spacer/page-break paragraphs immediately
# surrounding it, was removed from the Requisitos section.
$findRange = $d.Content
$found = $findRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Target paragraph 'Ver no Jupiter Salvar em pdf Salvar em docx' not found"
}

$target = $findRange.Paragraphs(1)

# Layout around the target paragraph is:
#   [blank]  [[Ver no Jupiter ...]]  [blank]  [blank + pageBreakBefore]  [blank]  [blank + pageBreakBefore]
# The first four of those six paragraphs (the blank one before the target,
# the target itself, and the two paragraphs following it) are deleted as a
# block; the trailing two paragraphs stay untouched.
$prev = $target.Previous()
$afterTarget = $target.Next()
$afterAfter = $afterTarget.Next()

$startPos = $prev.Range.Start
$endPos = $afterAfter.Range.End

$r = $d.Range($startPos, $endPos)
$r.Delete()
